$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.638.14"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.645.16"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "324.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0815"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "3.059.45"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "2.637.40"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.862"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "49.543.58"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0814"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "126.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").Value = "2.064.81"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.93%  "
